$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New master-location rows appended as part of the "2nd May Data Refresh".
# Columns: code, name, hierarchy_level, hierarchy_level_name, parent_loc_code,
#          lang_code, is_active, cr_by, cr_dtimes
$newRows = @(
    @(10113, 10113, 5, "الرمز البريدي", "BNMR", "ara", $true, "superadmin", "now()"),
    @(10114, 10114, 5, "الرمز البريدي", "BNMR", "ara", $true, "superadmin", "now()")
)

$r = 120
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
    $r = $r + 1
}

# Reflect the post-paste selection: whole rows from 122 to the sheet bottom.
$ws.Range("A122:A1048576").EntireRow.Select()
